$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data rows I2:J77
$data = @{
    2 = @(7, 8)
    3 = @(8, 8)
    4 = @(7, 7)
    5 = @(8, 8)
    6 = @(8, 8)
    7 = @(8, 8)
    8 = @(8, 8)
    9 = @(8, 8)
    10 = @(11, 11)
    11 = @(7, 7)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(8, 8)
    20 = @(8, 8)
    21 = @(8, 8)
    22 = @(9, 9)
    23 = @(8, 8)
    24 = @(8, 8)
    25 = @(8, 8)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(8, 8)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(8, 8)
    32 = @(8, 8)
    33 = @(9, 9)
    34 = @(8, 8)
    35 = @(8, 8)
    36 = @(9, 9)
    37 = @(8, 8)
    38 = @(8, 8)
    39 = @(7, 8)
    40 = @(8, 8)
    41 = @(8, 8)
    42 = @(8, 8)
    43 = @(10, 10)
    44 = @(7, 7)
    45 = @(7, 7)
    46 = @(7, 7)
    47 = @(7, 7)
    48 = @(7, 7)
    49 = @(6, 6)
    50 = @(6, 6)
    51 = @(7, 7)
    52 = @(7, 7)
    53 = @(5, 6)
    54 = @(7, 7)
    55 = @(7, 7)
    56 = @(6, 7)
    57 = @(6, 7)
    58 = @(7, 7)
    59 = @(7, 7)
    60 = @(7, 7)
    61 = @(7, 7)
    62 = @(7, 8)
    63 = @(6, 6)
    64 = @(7, 7)
    65 = @(9, 9)
    66 = @(6, 6)
    67 = @(10, 10)
    68 = @(9, 9)
    69 = @(9, 9)
    70 = @(6, 6)
    71 = @(6, 7)
    72 = @(4, 4)
    73 = @(8, 8)
    74 = @(6, 6)
    75 = @(6, 6)
    76 = @(7, 7)
    77 = @(5, 5)
}

foreach ($r in $data.Keys) {
    $pair = $data[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
